$d = $word.ActiveDocument

# The document's last paragraph currently holds "可预见性" with a pPr that
# carries Yu Mincho / ja-JP paragraph-mark formatting. The target state
# splits that single paragraph into three: a plain "可预见性" paragraph,
# a plain "等离子喷涂：" paragraph, and a long descriptive paragraph (with
# a simplified pPr, keeping only rFonts hint="eastAsia") made of many runs.

$lastIndex = $d.Paragraphs.Count
$target = $d.Paragraphs($lastIndex)

$targetStart = $target.Range.Start
$targetEnd = $target.Range.End

# Make sure the paragraph we are about to rewrite via InsertXML is not the
# very last thing in the document body (inserting a multi-paragraph XML
# fragment into a range that reaches the document's absolute end corrupts
# earlier content in this host) -- add a throwaway paragraph after it first.
$target.Range.InsertParagraphAfter()

# Range spanning the whole original paragraph, including its own paragraph
# mark, so the replacement XML's own <w:pPr> (or lack thereof) takes over.
$fullRange = $d.Range($targetStart, $targetEnd)

$xml = @'
<?xml version="1.0" standalone="yes"?>
<?mso-application progid="Word.Document"?>
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">
<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">
<pkg:xmlData>
<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
<w:body>
<w:p><w:r><w:rPr><w:rFonts w:hint="eastAsia"/></w:rPr><w:t>可预见性</w:t></w:r></w:p>
<w:p><w:r><w:rPr><w:rFonts w:hint="eastAsia"/></w:rPr><w:t>等离子喷涂：</w:t></w:r></w:p>
<w:p><w:pPr><w:rPr><w:rFonts w:hint="eastAsia"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:hint="eastAsia"/></w:rPr><w:t>等离子喷涂是以等离子弧作热源</w:t></w:r><w:r><w:t>将</w:t></w:r><w:r><w:rPr><w:rFonts w:hint="eastAsia"/></w:rPr><w:t>喷涂粉末材料在</w:t></w:r><w:proofErr w:type="gramStart"/><w:r><w:rPr><w:rFonts w:hint="eastAsia"/></w:rPr><w:t>等离子弧焰流中</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:rPr><w:rFonts w:hint="eastAsia"/></w:rPr><w:t>加热到熔化或半熔化状态</w:t></w:r><w:r><w:t>用高速气流将其吹成微小</w:t></w:r><w:r><w:rPr><w:rFonts w:hint="eastAsia"/></w:rPr><w:t>颗粒</w:t></w:r><w:r><w:rPr><w:rFonts w:hint="eastAsia"/></w:rPr><w:t>，</w:t></w:r><w:r><w:rPr><w:rFonts w:hint="eastAsia"/></w:rPr><w:t>喷射到经过处理的工件表面</w:t></w:r><w:r><w:rPr><w:rFonts w:hint="eastAsia"/></w:rPr><w:t>，</w:t></w:r><w:r><w:rPr><w:rFonts w:hint="eastAsia"/></w:rPr><w:t>形成牢固的覆盖层</w:t></w:r></w:p>
</w:body>
</w:document>
</pkg:xmlData>
</pkg:part>
</pkg:package>
'@

$fullRange.InsertXML($xml)

# Drop the throwaway trailing paragraph (it inherited the old Yu Mincho /
# ja-JP paragraph-mark formatting we no longer want).
$trailing = $d.Paragraphs($d.Paragraphs.Count)
$trailing.Range.Delete()
